$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new columns I/J on existing rows (per-month averages) ---
$ws.Range("I2").Formula = "=H2/6"
$ws.Range("J2").Formula = "=95/6"
$ws.Range("I7").Formula = "=172/6"
$ws.Range("J7").Formula = "=71/8"

# --- column widths for B:C (used by new table at rows 13-15) ---
$ws.Range("B:C").ColumnWidth = 17.85

# --- Row 13: header row for second table, wrapped + italic style ---
$ws.Range("A13").Value = "Indicator"
$ws.Range("A13").Font.Bold = $true

# (establish shared-string order: Homes visited(18) must precede the "Number/Total" strings)
$ws.Range("A14").Value = "Homes visited"

$ws.Range("B13").Value = "Number per month (1/2013-6/2013)"
$ws.Range("C13").Value = "Number per month (7/2013-3/2014)"
$ws.Range("F13").Value = "Total number (1/2013-6/2013)"
$ws.Range("G13").Value = "Total number (7/2013-3/2014)"
$ws.Range("B13:C13,F13:G13").Font.Italic = $true
$ws.Range("B13:C13,F13:G13").WrapText = $true
$ws.Range("A13:G13").RowHeight = 60

# --- Row 14: Homes visited ---
$ws.Range("B14").Value = 236
$ws.Range("C14").Formula = "=G14/8"
$ws.Range("C14").NumberFormat = "0"
$ws.Range("F14").Value = 1414
$ws.Range("G14").Value = 95

# --- Row 15: Births ---
$ws.Range("A15").Value = "Births"
$ws.Range("B15").Formula = "=F15/6"
$ws.Range("B15").NumberFormat = "0"
$ws.Range("C15").Formula = "=G15/8"
$ws.Range("C15").NumberFormat = "0"
$ws.Range("F15").Value = 172
$ws.Range("G15").Value = 71

# --- Row 19: third table header (indicator / homes visited / births) ---
$ws.Range("A19").Value = "Indicator"
$ws.Range("A19").Font.Bold = $true
$ws.Range("B19").Value = "Homes visited"
$ws.Range("C19").Value = "Births"

# --- Row 20: Total number (1/2013-6/2013) ---
$ws.Range("A20").Value = "Total number (1/2013-6/2013)"
$ws.Range("A20").Font.Italic = $true
$ws.Range("B20").Value = 1414
$ws.Range("C20").Value = 172

# --- Row 21: Total number (7/2013-3/2014) ---
$ws.Range("A21").Value = "Total number (7/2013-3/2014)"
$ws.Range("A21").Font.Italic = $true
$ws.Range("B21").Value = 95
$ws.Range("C21").Value = 71

# --- Row 22: Number per month (1/2013-6/2013) ---
$ws.Range("A22").Value = "Number per month (1/2013-6/2013)"
$ws.Range("A22").Font.Italic = $true
$ws.Range("B22").Value = 236
$ws.Range("C22").Formula = "=C20/6"
$ws.Range("C22").NumberFormat = "0"

# --- Row 23: Number per month (7/2013-3/2014) ---
$ws.Range("A23").Value = "Number per month (7/2013-3/2014)"
$ws.Range("A23").Font.Italic = $true
$ws.Range("B23").Formula = "=B21/8"
$ws.Range("B23").NumberFormat = "0"
$ws.Range("C23").Formula = "=C21/8"
$ws.Range("C23").NumberFormat = "0"

# --- selection matches target ---
$ws.Range("A13:C15").Select()
